# Insert a new weekly price record for "Pepino dulce" (Vega Modelo de Temuco)
# as row 94, pushing every subsequent record down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 94, shifting rows 94:145 down to 95:146.
$ws.Range("A94").EntireRow.Insert()

# Populate the newly inserted row 94 with the new observation.
$ws.Range("A94").Value = 10
$ws.Range("B94").Value = "Vega Modelo de Temuco"
$ws.Range("C94").Value = "La Araucanía"
$ws.Range("D94").Value = 44438
$ws.Range("E94").Value = 9
$ws.Range("F94").Value = 100112043
$ws.Range("G94").Value = "Pepino dulce"
$ws.Range("H94").Value = "Cultivar IV Región"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 160
$ws.Range("K94").Value = 21000
$ws.Range("L94").Value = 22000
$ws.Range("M94").Value = 21500
$ws.Range("N94").Value = "`$/bandeja 18 kilos"
$ws.Range("O94").Value = "Provincia de Limarí"
$ws.Range("P94").Value = 1194
$ws.Range("Q94").Value = 18
$ws.Range("R94").Value = "Hortaliza"
